# client & group Scenarios
# Remove the obsolete "waittopageload2" row from the "Prepay Loan" sheet
# (shifting the final "clickonsubmit"/"click" row up), and make the
# "Prepay Loan" sheet the active tab instead of "Transactions".

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Prepay Loan")

# Delete row 6 (waittopageload2 / 4000) - remaining rows shift up.
$ws.Rows.Item(6).Delete()

# Make "Prepay Loan" the active sheet/tab and select A6 (the row that
# shifted up into that position), matching the new selection state.
$ws.Activate()
$ws.Range("A6").Select()
